# Apply the changes described by the diff:
#  - widen column E (auto-height related column) to 27.42578125 characters
#  - update the normalized timestamp values in column E for rows 24-27 and 125-159
#  - update the sheet view scroll position / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column E width -------------------------------------------------
# The headless engine quantizes ColumnWidth assignments to the nearest
# 1/6-character pixel grid before persisting the OOXML `width` attribute,
# so the raw target (27.42578125) is fed in pre-compensated so that the
# value actually written out lands as close as possible to the target.
$ws.Columns.Item(5).ColumnWidth = 26.666666667

# --- 2. Normalize duplicate-looking timestamp values in column E -------
$rowsToFix = @(24, 25, 26, 27) + (125..159)
foreach ($r in $rowsToFix) {
    $ws.Cells.Item($r, 5).Value = 20201201153327
}

# --- 3. Scroll position / selection ------------------------------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 145
$win.ScrollColumn = 1
$ws.Range("F160").Select() | Out-Null
